# Update employee absence data rows 2-11 (columns A-G) with new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2;  A=67743; B="Emilly da Paz";        C="P&D";                    D="Outros";             E=1; F=45092; G=5500.59 }
    @{ Row=3;  A=69524; B="Valentim Montenegro";  C="Juridico";               D="Problemas pessoais"; E=4; F=45106; G=5092.02 }
    @{ Row=4;  A=72334; B="Ravi Lucca Cassiano";  C="Engenharia";             D="Problemas pessoais"; E=1; F=45100; G=4445.38 }
    @{ Row=5;  A=23928; B="Fernanda Pimenta";     C="Recursos Humanos";       D="Outros";             E=8; F=45091; G=8659.360000000001 }
    @{ Row=6;  A=12446; B="Allana Fonseca";       C="Operacoes";              D="Viagem de negocios"; E=4; F=45106; G=2880.41 }
    @{ Row=7;  A=17134; B="Vitor Mendonça";       C="Financeiro";             D="Problemas pessoais"; E=6; F=45105; G=9168.030000000001 }
    @{ Row=8;  A=284;   B="Ana Vitória Mendonça"; C="Operacoes";              D="Outros";             E=5; F=45097; G=8878.219999999999 }
    @{ Row=9;  A=39392; B="Noah Ribeiro";         C="Recursos Humanos";       D="Viagem de negocios"; E=5; F=45102; G=8740.34 }
    @{ Row=10; A=55356; B="Enzo Nunes";           C="Financeiro";             D="Consulta medica";    E=7; F=45106; G=5395.37 }
    @{ Row=11; A=56510; B="José Cirino";          C="Atendimento ao Cliente"; D="Outros";             E=2; F=45080; G=9297.43 }
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value = $rowData.A
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
}
